$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Replace the "Mark Hakkarinen" block (heading + "Outreach Ambassador"
#    + 'Editor of "..."' + "Email") with a new "LilyDaVine" heading
#    followed by a single "Outreach Support" paragraph.
# ---------------------------------------------------------------------
$startIdx = Get-ParaIndexByText $d "Mark Hakkarinen"
$endIdx   = Get-ParaIndexByText $d "Email"

$startPara = $d.Paragraphs($startIdx)
$endPara   = $d.Paragraphs($endIdx)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$block1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="105" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr><w:t>LilyDaVine</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="3B3B3B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Outreach Support</w:t></w:r></w:p>'

$rng.InsertXML($block1)

# ---------------------------------------------------------------------
# 2) Delete the "Outreach Support / Carlos Santiago / Outreach Support /
#    Emilio" run of paragraphs that used to follow "Semptly".
# ---------------------------------------------------------------------
$semptlyIdx = Get-ParaIndexByText $d "Semptly"
$p1 = $d.Paragraphs($semptlyIdx + 1)
$p2 = $d.Paragraphs($semptlyIdx + 4)
$delRng = $d.Range($p1.Range.Start, $p2.Range.End)
$delRng.Delete()

# ---------------------------------------------------------------------
# 3) Delete the "emelia / auditor / Nitego / Senior QA Tester" paragraphs
#    that follow "Release Coordinator".
# ---------------------------------------------------------------------
$relIdx = Get-ParaIndexByText $d "Release Coordinator"
$p3 = $d.Paragraphs($relIdx + 1)
$p4 = $d.Paragraphs($relIdx + 4)
$delRng2 = $d.Range($p3.Range.Start, $p4.Range.End)
$delRng2.Delete()
